$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (Sheet1 -> Translations)
$ws.Name = "Translations"

# --- Row 1 (headers): insert new "Entity Id" column, shift Type/Index right, drop "Id" ---
$ws.Range("A1").Value = "Entity Id"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Index"
# D1 "Original" and E1 "Translation" stay as-is

# --- Row 2 ---
$ws.Range("A2").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B2").Value = "Title"
$ws.Range("C2").ClearContents()
# D2 "Orig" and E2 "title" stay as-is

# --- Row 3 ---
$ws.Range("A3").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B3").Value = "ValidationMessage"
$ws.Range("C3").Value = 1
# D3 "Orig" and E3 "validation message" stay as-is

# --- Row 4 ---
$ws.Range("A4").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B4").Value = "Instruction"
$ws.Range("C4").ClearContents()
# D4 "Orig" and E4 "instruction" stay as-is

# --- Row 5 ---
$ws.Range("A5").Value = "AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$ws.Range("B5").Value = "OptionTitle"
$ws.Range("C5").Value = 2
# D5 "Orig" and E5 "option" stay as-is

# Column widths for the newly populated A/B/C columns (closest reachable values;
# internal width-quantization of the runtime rounds to 1/6-character steps)
$ws.Columns("A").ColumnWidth = 42.42
$ws.Columns("B").ColumnWidth = 17.25
$ws.Columns("C").ColumnWidth = 5.1

# Update the active cell / selection
$ws.Range("E10").Select() | Out-Null
